$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.448.39"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.569.05"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "'290.02"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "'0.3688"
$ws.Range("E7").Value = "  -1.45%  "
$ws.Range("D8").Value = "'49.92"
$ws.Range("E8").Value = "  +1.72%  "
$ws.Range("D9").Value = "'0.3375"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "'1.146"
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("D11").Value = "'0.07529"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "'21.14"
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("D14").Value = "'6.028"
$ws.Range("E14").Value = "  +2.03%  "
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").Value = "1.570.90"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "'90.35"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").Value = "'0.06773"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").Value = "'6.358"
$ws.Range("E21").Value = "  +3.20%  "
$ws.Range("D22").Value = "'16.39"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "'12.20"
$ws.Range("E23").Value = "  +3.05%  "
$ws.Range("D24").Value = "22.454.97"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").Value = "'2.372"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").Value = "'2.651"
$ws.Range("E26").Value = "  -1.88%  "
$ws.Range("D27").Value = "'20.01"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("D29").Value = "'5.055"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").Value = "'124.92"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "1.749.94"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").Value = "'1.068"
$ws.Range("E32").Value = "  +9.13%  "
$ws.Range("D33").Value = "'6.197"
$ws.Range("E33").Value = "  +4.13%  "
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").Value = "'9.793"
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("D36").Value = "'0.08339"
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").Value = "'0.02468"
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("D38").Value = "'1.357"
$ws.Range("E38").Value = "  -3.67%  "
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("D40").Value = "'0.06552"
$ws.Range("E40").Value = "  +2.07%  "
$ws.Range("D41").Value = "'5.402"
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("D42").Value = "'11.21"
$ws.Range("E42").Value = "  +2.36%  "
$ws.Range("D43").Value = "'0.6225"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'14.12"
$ws.Range("E44").Value = "  +2.22%  "
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").Value = "'3.801"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").Value = "'0.5856"
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("D48").Value = "'2.067"
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("D49").Value = "'127.72"
$ws.Range("E49").Value = "  +3.06%  "
$ws.Range("D50").Value = "'1.245"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").Value = "'0.07302"
$ws.Range("E51").Value = "  -0.13%  "
